# Update evaluation metrics across the three worksheets to reflect the
# final evaluation results for isolation_forest/augmented/none_4/split_1/test_50_50

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Summary": single row of aggregate metrics
# Columns: A=Model B=Accuracy C=Precision D=Recall E=F1 F=F2 G=F5
#          H=AUC I=TP J=FP K=TN L=FN
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.6226591760299626
$wsSummary.Range("C2").Value = 0.5771495877502945
$wsSummary.Range("D2").Value = 0.9176029962546817
$wsSummary.Range("E2").Value = 0.7086044830079538
$wsSummary.Range("F2").Value = 0.8207705192629816
$wsSummary.Range("G2").Value = 0.8972462849496443
$wsSummary.Range("H2").Value = 0.7754176661195977
$wsSummary.Range("I2").Value = 490
$wsSummary.Range("J2").Value = 359
$wsSummary.Range("K2").Value = 175
$wsSummary.Range("L2").Value = 44

# ---------------------------------------------------------------------
# Sheet "Classification Report": per-class precision/recall/f1/support
# ---------------------------------------------------------------------
$wsReport = $wb.Worksheets.Item("Classification Report")

# Row 2 - class "0"
$wsReport.Range("B2").Value = 0.7990867579908676
$wsReport.Range("C2").Value = 0.3277153558052435
$wsReport.Range("D2").Value = 0.4648074369189907

# Row 3 - class "1"
$wsReport.Range("B3").Value = 0.5771495877502945
$wsReport.Range("C3").Value = 0.9176029962546817
$wsReport.Range("D3").Value = 0.7086044830079538

# Row 4 - accuracy
$wsReport.Range("B4").Value = 0.6226591760299626
$wsReport.Range("C4").Value = 0.6226591760299626
$wsReport.Range("D4").Value = 0.6226591760299626
$wsReport.Range("E4").Value = 0.6226591760299626

# Row 5 - macro avg
$wsReport.Range("B5").Value = 0.688118172870581
$wsReport.Range("C5").Value = 0.6226591760299626
$wsReport.Range("D5").Value = 0.5867059599634722

# Row 6 - weighted avg
$wsReport.Range("B6").Value = 0.6881181728705811
$wsReport.Range("C6").Value = 0.6226591760299626
$wsReport.Range("D6").Value = 0.5867059599634722

# ---------------------------------------------------------------------
# Sheet "Confusion Matrix"
# ---------------------------------------------------------------------
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 - Actual 0
$wsConfusion.Range("B2").Value = 175
$wsConfusion.Range("C2").Value = 359

# Row 3 - Actual 1
$wsConfusion.Range("B3").Value = 44
$wsConfusion.Range("C3").Value = 490
